# Apply the cryptos list update (values taken from the upstream commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.904.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.69'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.74%  '
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0615'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.871.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.638.39'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.575'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.17'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.894.85'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.13'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.38'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.111'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0484'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.13'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.52%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.420.41'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('E37').Value = '  +2.39%  '
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.920'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.556'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('E45').Value = '  +3.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.39%  '
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.780.57'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.100'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('E51').Value = '  +0.68%  '
